# The workbook lists required R packages, one per row, in column A of Sheet1.
# This edit appends a new package name ("shinythemes") to the end of that list,
# which is the content change captured by the diff (new shared string +
# new row 52, dimension growing from A1:A51 to A1:A52).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A, then target the row right after it.
$lastRow = $ws.Range("A1048576").End(-4162).Row
$newRow = $lastRow + 1

# Write the new package name into the first empty cell of the list.
$ws.Cells.Item($newRow, 1).Value = "shinythemes"

# Select the newly added cell and scroll it into view, matching the
# updated sheet view (active cell / visible range) recorded after the edit.
$ws.Range("A" + $newRow).Select()
$excel.ActiveWindow.ScrollRow = [Math]::Max(1, $newRow - 3)
